$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 508 (pushes the existing row 508..572 down to 509..573)
$ws.Rows("508").Insert()

# Populate the newly inserted row with the new weekly price record
$ws.Cells.Item(508, 1).Value = 4
$ws.Cells.Item(508, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(508, 3).Value = "Los Lagos"
$ws.Cells.Item(508, 4).Value = 45131
$ws.Cells.Item(508, 5).Value = 10
$ws.Cells.Item(508, 6).Value = 100112023
$ws.Cells.Item(508, 7).Value = "Brócoli"
$ws.Cells.Item(508, 8).Value = "Sin especificar"
$ws.Cells.Item(508, 9).Value = "Primera"
$ws.Cells.Item(508, 10).Value = 500
$ws.Cells.Item(508, 11).Value = 1500
$ws.Cells.Item(508, 12).Value = 1500
$ws.Cells.Item(508, 13).Value = 1500
$ws.Cells.Item(508, 14).Value = "$/unidad"
$ws.Cells.Item(508, 15).Value = "Región Metropolitana"
$ws.Cells.Item(508, 16).Value = 1500
$ws.Cells.Item(508, 17).Value = 1
$ws.Cells.Item(508, 18).Value = "Hortaliza"
